# BYOC Motherboard update for ASUS
# Populate the "华硕主板规格" (ASUS motherboard spec) sheet with a
# 3-motherboard comparison table (B650M-K / B650M-R / B650M-AYW).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("华硕主板规格")

# Column widths: A narrower (label column), B:D wider & equal (data columns)
$ws.Columns("A:A").ColumnWidth = 14.125
$ws.Columns("B:D").ColumnWidth = 31.625

$data = @(
    @("华硕主板",   "B650M-K",                      "B650M-R",                     "B650M-AYW"),
    @("CPU插槽",    "AM5",                          "AM5",                         "AM5"),
    @("内存插槽",   "DDR5x2",                       "DDR5x2",                      "DDR5x2"),
    @("扩展插槽",   "1个 PCIe 4.0 x16",             "1个 PCIe 4.0 x16",            "2个 PCIe 4.0 x16"),
    @("",           "2个 PCIe 4.0 x1",              "1个 PCIe 4.0 x1",             ""),
    @("M.2插槽",    "2个个 M.2 PCIe 4.0 x4",        "2个个 M.2 PCIe 4.0 x4",       "2个个 M.2 PCIe 5.0 x4 + 4.0 x4"),
    @("SATA",       "4个 SATA 6Gb/s",               "4个 SATA 6Gb/s",              "4个 SATA 6Gb/s"),
    @("显示接口",   "1个 HDMI",                     "1个 HDMI",                    "1个 HDMI"),
    @("",           "1个 VGA",                      "",                            ""),
    @("USB接口（后）", "1个 USB 2.0 (BIOS FlashBack)", "",                         ""),
    @("",           "3个 USB 2.0",                  "4个 USB 2.0",                 ""),
    @("",           "2个 USB 3.2 Gen1",             "",                            ""),
    @("",           "2个 USB 3.2 Gen2 Type-A",      "2个 USB 5Gbps",               ""),
    @("Wi-Fi",      "无",                           "无",                          "Wi-Fi 6 + 蓝牙5.3")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 1
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Length; $j++) {
        $val = $rowVals[$j]
        if ($val -ne "") {
            $ws.Cells.Item($r, $j + 1).Value = $val
        }
    }
}
